$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 120
$ws.Cells.Item(33, 9).Value = 130.58333
$ws.Cells.Item(33, 11).Value = 130.58333
$ws.Cells.Item(33, 13).Value = 98.41667000000001

$ws.Cells.Item(70, 8).Value = 2900
$ws.Cells.Item(70, 10).Value = 3000
$ws.Cells.Item(70, 12).Value = 9000
$ws.Cells.Item(70, 14).Value = -9540

$ws.Cells.Item(73, 8).Value = 2900
$ws.Cells.Item(73, 10).Value = 3000
$ws.Cells.Item(73, 12).Value = 9000
$ws.Cells.Item(73, 14).Value = -10872

$ws.Cells.Item(74, 8).Value = 1566.6666
$ws.Cells.Item(74, 9).Value = 1566.6666
$ws.Cells.Item(74, 11).Value = 1566.6666
$ws.Cells.Item(74, 13).Value = -630.6666

$ws.Cells.Item(77, 8).Value = 1566.6666
$ws.Cells.Item(77, 9).Value = 1566.6666
$ws.Cells.Item(77, 11).Value = 7833.333000000001
$ws.Cells.Item(77, 13).Value = -3153.333000000001

$ws.Cells.Item(98, 8).Value = 6153.3335
$ws.Cells.Item(98, 9).Value = 6153.3335
$ws.Cells.Item(98, 11).Value = 6153.3335
$ws.Cells.Item(98, 13).Value = -4655.3335

$ws.Cells.Item(111, 8).Value = 3439.2727
$ws.Cells.Item(111, 9).Value = 266.8
$ws.Cells.Item(111, 10).Value = 6083
$ws.Cells.Item(111, 11).Value = 800.4000000000001
$ws.Cells.Item(111, 12).Value = 18249
$ws.Cells.Item(111, 13).Value = 2266.6
$ws.Cells.Item(111, 14).Value = -24383

$ws.Cells.Item(122, 8).Value = 6153.3335
$ws.Cells.Item(122, 9).Value = 6153.3335
$ws.Cells.Item(122, 11).Value = 18460.0005
$ws.Cells.Item(122, 13).Value = -16010.0005

$ws.Cells.Item(131, 8).Value = 392.5
$ws.Cells.Item(131, 9).Value = 392.5
$ws.Cells.Item(131, 11).Value = 1177.5
$ws.Cells.Item(131, 13).Value = 3862.5

$ws.Cells.Item(137, 8).Value = 1172.625
$ws.Cells.Item(137, 9).Value = 1275.7858
$ws.Cells.Item(137, 10).Value = 450.5
$ws.Cells.Item(137, 11).Value = 3827.3574
$ws.Cells.Item(137, 12).Value = 1351.5
$ws.Cells.Item(137, 13).Value = -1277.3574
$ws.Cells.Item(137, 14).Value = -6451.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6748.7144
$ws.Cells.Item(32, 9).Value = 5334
$ws.Cells.Item(32, 10).Value = 11275.8
$ws.Cells.Item(32, 11).Value = 5334
$ws.Cells.Item(32, 12).Value = 11275.8
$ws.Cells.Item(32, 13).Value = -5047
$ws.Cells.Item(32, 14).Value = -11849.8

$ws.Cells.Item(43, 8).Value = 39999
$ws.Cells.Item(43, 10).Value = 39999
$ws.Cells.Item(43, 12).Value = 39999
$ws.Cells.Item(43, 14).Value = -40625

$ws.Cells.Item(61, 8).Value = 1999.25
$ws.Cells.Item(61, 9).Value = 1999.25
$ws.Cells.Item(61, 11).Value = 1999.25
$ws.Cells.Item(61, 13).Value = -1787.25

$ws.Cells.Item(110, 8).Value = 5713.2856
$ws.Cells.Item(110, 10).Value = 5999.6
$ws.Cells.Item(110, 12).Value = 5999.6
$ws.Cells.Item(110, 14).Value = -10089.6

$ws.Cells.Item(122, 8).Value = 1974.9166
$ws.Cells.Item(122, 9).Value = 1699.909
$ws.Cells.Item(122, 11).Value = 5099.727000000001
$ws.Cells.Item(122, 13).Value = -2649.727000000001

$ws.Cells.Item(132, 8).Value = 2285.7058
$ws.Cells.Item(132, 9).Value = 1928.5625
$ws.Cells.Item(132, 11).Value = 5785.6875
$ws.Cells.Item(132, 13).Value = -3255.6875

$ws.Cells.Item(136, 8).Value = 1999.25
$ws.Cells.Item(136, 9).Value = 1999.25
$ws.Cells.Item(136, 11).Value = 5997.75
$ws.Cells.Item(136, 13).Value = -3447.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3452.647
$ws.Cells.Item(86, 9).Value = 1979.5
$ws.Cells.Item(86, 11).Value = 1979.5
$ws.Cells.Item(86, 13).Value = -856.5

$ws.Cells.Item(89, 8).Value = 3452.647
$ws.Cells.Item(89, 9).Value = 1979.5
$ws.Cells.Item(89, 11).Value = 9897.5
$ws.Cells.Item(89, 13).Value = -4281.5

$ws.Cells.Item(134, 8).Value = 2846
$ws.Cells.Item(134, 9).Value = 2846
$ws.Cells.Item(134, 11).Value = 8538
$ws.Cells.Item(134, 13).Value = -6003

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7554.875
$ws.Cells.Item(31, 9).Value = 7902.647
$ws.Cells.Item(31, 10).Value = 6710.2856
$ws.Cells.Item(31, 11).Value = 7902.647
$ws.Cells.Item(31, 12).Value = 6710.2856
$ws.Cells.Item(31, 13).Value = -7607.647
$ws.Cells.Item(31, 14).Value = -7300.2856

$ws.Cells.Item(34, 8).Value = 7554.875
$ws.Cells.Item(34, 9).Value = 7902.647
$ws.Cells.Item(34, 10).Value = 6710.2856
$ws.Cells.Item(34, 11).Value = 7902.647
$ws.Cells.Item(34, 12).Value = 6710.2856
$ws.Cells.Item(34, 13).Value = -7700.647
$ws.Cells.Item(34, 14).Value = -7114.2856

$ws.Cells.Item(39, 8).Value = 6313.722
$ws.Cells.Item(39, 9).Value = 1282.0714
$ws.Cells.Item(39, 11).Value = 1282.0714
$ws.Cells.Item(39, 13).Value = -891.0714

$ws.Cells.Item(49, 8).Value = 6313.722
$ws.Cells.Item(49, 9).Value = 1282.0714
$ws.Cells.Item(49, 11).Value = 1282.0714
$ws.Cells.Item(49, 13).Value = -1100.0714

$ws.Cells.Item(60, 8).Value = 24990
$ws.Cells.Item(60, 9).Value = 24990
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 24990
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = -24479
$ws.Cells.Item(60, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 535.3333
$ws.Cells.Item(50, 10).Value = 1506
$ws.Cells.Item(50, 12).Value = 4518
$ws.Cells.Item(50, 14).Value = -5480

$ws.Cells.Item(53, 8).Value = 535.3333
$ws.Cells.Item(53, 10).Value = 1506
$ws.Cells.Item(53, 12).Value = 4518
$ws.Cells.Item(53, 14).Value = -5480

$ws.Cells.Item(109, 8).Value = 370
$ws.Cells.Item(109, 9).Value = 370
$ws.Cells.Item(109, 11).Value = 1110
$ws.Cells.Item(109, 13).Value = -70

$ws.Cells.Item(131, 8).Value = 2639.889
$ws.Cells.Item(131, 9).Value = 3633.3333
$ws.Cells.Item(131, 10).Value = 2143.1667
$ws.Cells.Item(131, 11).Value = 10899.9999
$ws.Cells.Item(131, 12).Value = 6429.500100000001
$ws.Cells.Item(131, 13).Value = -5859.999899999999
$ws.Cells.Item(131, 14).Value = -16509.5001

$ws.Cells.Item(132, 8).Value = 2513
$ws.Cells.Item(132, 9).Value = 2432
$ws.Cells.Item(132, 11).Value = 21888
$ws.Cells.Item(132, 13).Value = -19358

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 32950
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 5945.8184
$ws.Cells.Item(132, 9).Value = 5740.5
$ws.Cells.Item(132, 11).Value = 17221.5
$ws.Cells.Item(132, 13).Value = -14691.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 887.8889
$ws.Cells.Item(46, 9).Value = 799.25
$ws.Cells.Item(46, 10).Value = 958.8
$ws.Cells.Item(46, 11).Value = 799.25
$ws.Cells.Item(46, 12).Value = 958.8
$ws.Cells.Item(46, 13).Value = -611.25
$ws.Cells.Item(46, 14).Value = -1334.8

$ws.Cells.Item(82, 8).Value = 1374.25
$ws.Cells.Item(82, 10).Value = 686.5
$ws.Cells.Item(82, 12).Value = 686.5
$ws.Cells.Item(82, 14).Value = -1408.5

$ws.Cells.Item(85, 8).Value = 1374.25
$ws.Cells.Item(85, 10).Value = 686.5
$ws.Cells.Item(85, 12).Value = 686.5
$ws.Cells.Item(85, 14).Value = -3182.5

$ws.Cells.Item(122, 8).Value = 10000
$ws.Cells.Item(122, 9).Value = 10000
$ws.Cells.Item(122, 11).Value = 30000
$ws.Cells.Item(122, 13).Value = -27550

$ws.Cells.Item(136, 8).Value = 7877.625
$ws.Cells.Item(136, 9).Value = 8003.6665
$ws.Cells.Item(136, 11).Value = 24010.9995
$ws.Cells.Item(136, 13).Value = -21460.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(55, 8).Value = 1047.5
$ws.Cells.Item(55, 9).Value = 1047.5
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 1047.5
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).Value = -770.5
$ws.Cells.Item(55, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 2960.923
$ws.Cells.Item(132, 9).Value = 3070.1667
$ws.Cells.Item(132, 11).Value = 9210.500100000001
$ws.Cells.Item(132, 13).Value = -6680.500100000001

$ws.Cells.Item(136, 8).Value = 4655.2856
$ws.Cells.Item(136, 9).Value = 4597.8335
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 13793.5005
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = -11243.5005
$ws.Cells.Item(136, 14).Value = -20100
